$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes old row 3 -> row 4)
$ws.Rows.Item(3).Insert()

# --- Row 2 (sCs / Nlgn1 / Nrxn2 / FAPs) : update numeric columns ---
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.2197353333333333
$ws.Cells.Item(2, 8).Value = 0.659206
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.039512
$ws.Cells.Item(2, 14).Value = 0.118536
$ws.Cells.Item(2, 15).Value = 0.006183359004302676
$ws.Cells.Item(2, 16).Value = 0.006183359004302677
$ws.Cells.Item(2, 17).Value = 0.008682182490666665
$ws.Cells.Item(2, 18).Value = 0.07813964241599999
$ws.Cells.Item(2, 19).Value = 0.006183359004302676
$ws.Cells.Item(2, 20).Value = 0.006183359004302677

# --- Row 3 (new): sCs / Nlgn1 / Nrxn2 / ECs ---
$ws.Cells.Item(3, 1).Value = "sCs"
$ws.Cells.Item(3, 2).Value = "Nlgn1"
$ws.Cells.Item(3, 3).Value = "Nrxn2"
$ws.Cells.Item(3, 4).Value = "ECs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.2197353333333333
$ws.Cells.Item(3, 8).Value = 0.659206
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 6.039054
$ws.Cells.Item(3, 14).Value = 18.117162
$ws.Cells.Item(3, 15).Value = 0.9450708374258476
$ws.Cells.Item(3, 16).Value = 0.9450708374258477
$ws.Cells.Item(3, 17).Value = 1.326993543708
$ws.Cells.Item(3, 18).Value = 11.942941893372
$ws.Cells.Item(3, 19).Value = 0.9450708374258476
$ws.Cells.Item(3, 20).Value = 0.9450708374258477

# --- Row 4 (was old row 3): sCs / Nlgn1 / Nrxn2 / sCs : update numeric columns ---
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.2197353333333333
$ws.Cells.Item(4, 8).Value = 0.659206
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.3114883333333333
$ws.Cells.Item(4, 14).Value = 0.934465
$ws.Cells.Item(4, 15).Value = 0.04874580356984966
$ws.Cells.Item(4, 16).Value = 0.04874580356984967
$ws.Cells.Item(4, 17).Value = 0.06844499275444443
$ws.Cells.Item(4, 18).Value = 0.61600493479
$ws.Cells.Item(4, 19).Value = 0.04874580356984966
$ws.Cells.Item(4, 20).Value = 0.04874580356984967
